# duncan2016.xlsx -- melt "cleaned" data to fit metafor analysis structure:
#   * average the two "constant" treatment rows (they don't share the same
#     mean as the flux treatments) into a new avg_constant / avg_constant_error
#     pair, highlighted so it's clear it's a derived summary row
#   * add a new (currently empty) "metafor calcs" sheet to hold the metafor
#     analysis that will consume this melted data
#   * leave "cleaned" as the active/selected sheet, now focused on the new rows

$wb = $excel.ActiveWorkbook

$cleaned = $wb.Worksheets.Item("cleaned")

# New sheet for the metafor calculations, placed right after "cleaned".
$metafor = $wb.Worksheets.Add($null, $cleaned)
$metafor.Name = "metafor calcs"

# Re-activate "cleaned" -- adding a sheet makes the new one active/selected.
$cleaned.Activate() | Out-Null

# Header row for the averaged-constant summary.
$cleaned.Range("B9").Value = "avg_constant"
$cleaned.Range("C9").Value = "avg_constant_error"

# Average the two "constant" treatment rows (rows 2 and 4) for y and avg_error.
$cleaned.Range("B10").Formula = "=AVERAGE(A2,A4)"
$cleaned.Range("C10").Formula = "=AVERAGE(C2,C4)"

# Highlight the new summary block in orange so it stands out as derived.
$cleaned.Range("B9:C10").Interior.Color = 49407

# Match the author's final selection/view state.
$cleaned.Range("B9:C10").Select() | Out-Null
